$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("reaction14")

# Clear out the existing row contents (C1:R1) so the sheet dimension shrinks back down
$ws.Range("A1:R1").ClearContents()

# Set the new values for the remaining two cells
$ws.Range("A1").Value = 28
$ws.Range("B1").Value = 29
